# Reorganize "Filed Cases" rows 7-10 and 21-23 (Foundation Filed Cases / Related
# Cases / Not Litigated updates), add two new case entries, and switch the
# active tab to "Not Litigated".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Filed Cases")

# Row 7
$ws.Range("A7").Value = 'California'
$ws.Range("B7").Value = '2000, 2001'
$ws.Range("C7").Value = 'Safe In Heaven Dead Productions, Sal Jenco, Johnny Depp'
$ws.Range("D7").Value = 'Anthony Fox, Family'
$ws.Range("E7").ClearContents()
$ws.Range("G7").Value = 'Fraud, Embellzement'
$ws.Range("I7").Value = 'Nightclubs'
$ws.Range("J7").Value = 'Nightclub Owner'
$ws.Range("K7").ClearContents()
$ws.Range("R7").Value = 'Fixed, Death, Dismissed'
$ws.Range("S7").Value = 'Paul Schindler'
$ws.Range("T7").Value = 'https://trellis.law/case/sc062176/anthony-v-fox-vs-safe-in-heaven-dead-productions-inc-et-al'
$ws.Range("U7").ClearContents()
$ws.Range("V7").ClearContents()

# Row 8
$ws.Range("A8").Value = 'Michigan'
$ws.Range("B8").Value = 2024
$ws.Range("C8").Value = 'Sean Combs, RICO'
$ws.Range("D8").Value = 'Derrick Lee Cardello-Smith'
$ws.Range("E8").Value = 'Properia Persona'
$ws.Range("G8").Value = 'Assault, Investments, Witnessing'
$ws.Range("I8").Value = 'Music, Properties, Investments'
$ws.Range("J8").Value = 'Music'
$ws.Range("K8").Value = 'Music, Business'
$ws.Range("R8").Value = '100,000,000 million default judgement'
$ws.Range("T8").Value = 'https://www.youtube.com/watch?v=OSP-z2qE914'
$ws.Range("U8").Value = 'https://www.youtube.com/watch?v=_rsV-w6_bKE'
$ws.Range("V8").Value = 'https://www.metrotimes.com/news/michigan-inmate-wins-100-million-judgment-against-sean-diddy-combs-for-sexual-assault-37301149'

# Row 9
$ws.Range("B9").Value = 2023
$ws.Range("C9").Value = 'David Miscavige, Church of Scientology International, Religious Technology Center'
$ws.Range("D9").Value = 'Leah Remini'
$ws.Range("E9").Value = 'Edwards Henderson Lehrman, EPLLC, Brittany N. Henderson, Bradley J. Edwards, Seth M. Lehrman, Motley Rise LLC, Linda Singer, Carmen S. Scott'
$ws.Range("G9").Value = 'Civil Harassment, Stalking, False Light, Defamation, Tortious Interference in Business'
$ws.Range("I9").ClearContents()
$ws.Range("J9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").Value = 'https://youtu.be/g2D70ddKnFc?si=_oKN8zLJ5MoihK_y'
$ws.Range("U9").Value = 'https://deadline.com/wp-content/uploads/2023/08/remini-scientology-amended-complaint.pdf'

# Row 10
$ws.Range("B10").Value = 2019
$ws.Range("C10").Value = 'Blackcube, David Boeis, Harvey Weinstein, Lisa Bloom'
$ws.Range("D10").Value = 'Rose McGowan'
$ws.Range("E10").ClearContents()
$ws.Range("F10").Value = 'David Boeis'
$ws.Range("G10").Value = 'Fraud, RICO, Banes Act'
$ws.Range("I10").Value = 'Film, Fixers, Investigators'
$ws.Range("T10").ClearContents()
$ws.Range("U10").ClearContents()

# Row 21
$ws.Range("A21").Value = 'California'
$ws.Range("B21").Value = 1997
$ws.Range("C21").Value = 'Bryan Singer'
$ws.Range("D21").Value = 'Ryan Glomboske, David Stockdale, St. Albin'
$ws.Range("E21").Value = 'Peter Gordon'
$ws.Range("G21").Value = 'Invasion of Privacy, Negligence, Infliction of Emotional Distress'
$ws.Range("I21").Value = 'Apt Pupil, Movie Set, Film'
$ws.Range("J21").Value = 'Child Actors'
$ws.Range("K21").Value = 'Director, Pheonix Pictures'
$ws.Range("U21").Value = 'https://ew.com/article/1997/05/02/indecent-proposal-set-apt-pupil/'

# Row 22
$ws.Range("A22").Value = 'Georgia'
$ws.Range("B22").ClearContents()
$ws.Range("C22").Value = 'Tyler Perry, Tyler Perry Studios, Brett Hendrix'
$ws.Range("D22").ClearContents()
$ws.Range("E22").ClearContents()
$ws.Range("G22").Value = 'Sexual Harassment Employment, Wrongful termination, Defamation'
$ws.Range("T22").ClearContents()
$ws.Range("U22").ClearContents()

# Row 23
$ws.Range("B23").Value = 2024
$ws.Range("C23").Value = 'Christian Combs, Sean Combs, '
$ws.Range("D23").Value = 'Grace O''Marcaigh'
$ws.Range("E23").Value = ' Tyrone Blackburn, Rodney S Digs, Tyrone S Ahmen, Ivie Mcneill Wyatt Purcell & Diggs, TA Blackburn Law PLLC'
$ws.Range("F23").ClearContents()
$ws.Range("G23").Value = 'Assault, Battery, Sexual Assault, Premises Liability, Aiding & Abetting, Intentional Infliction of Emotional Distress, Negligent Infliction of Emotional Distress, '
$ws.Range("I23").ClearContents()
$ws.Range("T23").Value = 'https://www.youtube.com/watch?v=JOPpS8i_voM'
$ws.Range("U23").Value = 'https://deadline.com/wp-content/uploads/2024/04/combs-conformed-suit.pdf'

# Switch active sheet from "Filed Cases" to "Not Litigated"
$wsNotLitigated = $wb.Worksheets.Item("Not Litigated")
$wsNotLitigated.Activate()
